# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
# (to preserve formats like trailing zeros / multi-dot price strings).
# Force text number format before assigning so Excel doesn't coerce them
# into real numbers.
$textCells = @(
    'D5',
    'D6',
    'D7',
    'D11',
    'D12',
    'D16',
    'D18',
    'D19',
    'D21',
    'D25',
    'D28',
    'D29',
    'D30',
    'D32',
    'D35',
    'D36',
    'D37',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '64.432.46'
$ws.Range('E2').Value = '  -3.39%  '
$ws.Range('D3').Value = '3.173.59'
$ws.Range('E3').Value = '  -5.04%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '571.75'
$ws.Range('E5').Value = '  -2.85%  '
$ws.Range('D6').Value = '169.37'
$ws.Range('E6').Value = '  -7.89%  '
$ws.Range('D7').Value = '0.606'
$ws.Range('E7').Value = '  -6.21%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '3.180.04'
$ws.Range('E9').Value = '  -4.73%  '
$ws.Range('E10').Value = '  -5.61%  '
$ws.Range('D11').Value = '6.81'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '0.391'
$ws.Range('E12').Value = '  -3.72%  '
$ws.Range('D13').Value = '3.722.21'
$ws.Range('E13').Value = '  -4.99%  '
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').Value = '64.502.58'
$ws.Range('E15').Value = '  -3.22%  '
$ws.Range('D16').Value = '25.40'
$ws.Range('E16').Value = '  -4.78%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.197.39'
$ws.Range('E17').Value = '  -3.28%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0000158'
$ws.Range('E18').Value = '  -5.05%  '
$ws.Range('D19').Value = '417.32'
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').Value = '5.33'
$ws.Range('E21').Value = '  -3.94%  '
$ws.Range('E22').Value = '  -3.51%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = '70.12'
$ws.Range('E25').Value = '  -2.90%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('D28').Value = '0.0000104'
$ws.Range('E28').Value = '  -9.89%  '
$ws.Range('D29').Value = '8.83'
$ws.Range('E29').Value = '  -2.54%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -5.26%  '
$ws.Range('D32').Value = '21.74'
$ws.Range('E32').Value = '  -3.40%  '
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('E34').Value = '  -3.71%  '
$ws.Range('D35').Value = '6.37'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('D36').Value = '1.13'
$ws.Range('E36').Value = '  -5.55%  '
$ws.Range('D37').Value = '156.25'
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('E38').Value = '  -5.96%  '
$ws.Range('E39').Value = '  -6.11%  '
$ws.Range('D40').Value = '2.710.37'
$ws.Range('E40').Value = '  -5.69%  '
$ws.Range('D41').Value = '4.23'
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('D42').Value = '24.27'
$ws.Range('E42').Value = '  -8.69%  '
$ws.Range('D43').Value = '39.07'
$ws.Range('D44').Value = '0.716'
$ws.Range('E44').Value = '  -6.10%  '
$ws.Range('D45').Value = '0.0623'
$ws.Range('E45').Value = '  -6.69%  '
$ws.Range('E46').Value = '  -8.44%  '
$ws.Range('E47').Value = '  -3.83%  '
$ws.Range('D48').Value = '21.54'
$ws.Range('E48').Value = '  -7.87%  '
$ws.Range('D49').Value = '292.09'
$ws.Range('E49').Value = '  -7.51%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.0992'
$ws.Range('E50').Value = '  -5.57%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  -0.14%  '
